$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.06328177979961902
$ws.Range("C2").Value = 0.05231270169004087
$ws.Range("D2").Value = 0.1529057820181812
$ws.Range("E2").Value = 0.4998867070740569
$ws.Range("G2").Value = 0.768386970581898
